$d = $word.ActiveDocument

# Locate the "KGB  connections" paragraph using Find, then grab its
# enclosing paragraph so we can rewrite the whole run sequence in one go.
$seek = $d.Content
$found = $seek.Find.Execute("KGB  connections")
if (-not $found) {
    Write-Output "KGB paragraph not found"
} else {
    $para = $seek.Paragraphs(1)
    $pRange = $para.Range

    # Exclude the trailing paragraph mark so InsertXML only replaces the
    # run content, leaving the paragraph (and its pPr) intact.
    $target = $d.Range($pRange.Start, $pRange.End - 1)

    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
      '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
      '<w:body>' +
      '<w:p>' +
        '<w:r><w:rPr><w:b/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">KGB  </w:t></w:r>' +
        '<w:r><w:rPr><w:b/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr><w:t>C</w:t></w:r>' +
        '<w:r><w:rPr><w:b/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>onnections</w:t></w:r>' +
        '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>: All face cards</w:t></w:r>' +
        '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr><w:t>and Aces</w:t></w:r>' +
        '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>gain Cripple 1 (</w:t></w:r>' +
        '<w:r><w:rPr><w:b/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Criple X</w:t></w:r>' +
        '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">: Opponent discards X random cards when played) </w:t></w:r>' +
        '<w:r><w:rPr><w:b/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr><w:t>IMPLEMENTED</w:t></w:r>' +
      '</w:p>' +
      '</w:body></w:document>' +
      '</pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($pkg)
    Write-Output "KGB Connections paragraph updated."
}
